$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.1
$ws.Range("H2").Value = 3.15
$ws.Range("I2").Value = 3.15
$ws.Range("J2").Value = 2.67
$ws.Range("K2").Value = 2.1
$ws.Range("L2").Value = 3.6
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 9
$ws.Range("O2").Value = 1.26
$ws.Range("P2").Value = 3.08
$ws.Range("Q2").Value = 1.91
$ws.Range("R2").Value = 1.7
$ws.Range("S2").Value = 1.37
$ws.Range("T2").Value = 2.5
$ws.Range("U2").Value = 1.77
$ws.Range("V2").Value = 2
$ws.Range("W2").Value = 6.2
$ws.Range("X2").Value = 8.25
$ws.Range("Y2").Value = 7.4
$ws.Range("Z2").Value = 16
$ws.Range("AA2").Value = 14.5
$ws.Range("AC2").Value = 9
$ws.Range("AD2").Value = 5.4
$ws.Range("AE2").Value = 11.5
$ws.Range("AF2").Value = 45
$ws.Range("AG2").Value = 300
$ws.Range("AH2").Value = 8.25
$ws.Range("AI2").Value = 14
$ws.Range("AJ2").Value = 9.25
$ws.Range("AK2").Value = 32
$ws.Range("AL2").Value = 22
$ws.Range("AM2").Value = 27
$ws.Range("AN2").Value = 4.05
$ws.Range("AO2").Value = 10.75
$ws.Range("AP2").Value = 18.5
$ws.Range("AQ2").Value = 40
$ws.Range("AR2").Value = 75
$ws.Range("AT2").Value = 2.55
$ws.Range("AU2").Value = 6.8
$ws.Range("AV2").Value = 55
$ws.Range("AW2").Value = 5.1
$ws.Range("AX2").Value = 17
$ws.Range("AY2").Value = 22
$ws.Range("AZ2").Value = 80
$ws.Range("BA2").Value = 110
$ws.Range("BB2").Value = 250
